$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the FileName, Author(s), Link (PDF), Edn/Vol, and Year values
# between row 36 (Friedberg/Insel/Spence) and row 37 (Hoffman/Kunze)
$cols = @("B", "D", "E", "F", "G")

foreach ($col in $cols) {
    $cell36 = $ws.Range($col + "36")
    $cell37 = $ws.Range($col + "37")
    $v36 = $cell36.Value2
    $v37 = $cell37.Value2
    $cell36.Value2 = $v37
    $cell37.Value2 = $v36
}
